# correção nos dados e inicio da analise PNAD 2009
#
# The source DataFrame dropped the empty "grandes regiões e unidades da
# federação" section header row (which had no numeric data associated
# with it). Re-exporting the corrected data shifts every data row that
# followed up by one row and removes the now-unused trailing row.
#
# Deleting worksheet row 6 (the row that only contained the
# "grandes regiões e unidades da federação" label in column A, with no
# values in B:G) reproduces exactly that: rows 7:37 shift up to become
# rows 6:36, the sheet dimension shrinks from A1:G37 to A1:G36, and the
# now-unused shared string for "grandes regiões e unidades da federação"
# is dropped from the workbook's shared string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:6").Delete()
